# Update cryptos.xlsx crypto listing with latest prices / % changes
# as scraped on Mon Apr  1 20:38:42 UTC 2024 via GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "69.442.14"
$ws.Cells.Item(2, 5).Value = "  -1.84%  "
$ws.Cells.Item(3, 4).Value = "3.478.05"
$ws.Cells.Item(3, 5).Value = "  -4.22%  "
$ws.Cells.Item(4, 5).Value = "  +0.13%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "579.34"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -4.22%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "192.47"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -3.28%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.612"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -2.51%  "
$ws.Cells.Item(8, 4).Value = "3.466.60"
$ws.Cells.Item(8, 5).Value = "  -4.22%  "
$ws.Cells.Item(9, 5).Value = "  -0.01%  "
$ws.Cells.Item(10, 5).Value = "  -7.81%  "
$ws.Cells.Item(11, 5).Value = "  -4.64%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "51.36"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -4.65%  "
$ws.Cells.Item(13, 5).Value = "  -6.74%  "
$ws.Cells.Item(14, 5).Value = "  -4.40%  "
$ws.Cells.Item(15, 4).Value = "4.034.44"
$ws.Cells.Item(15, 5).Value = "  -4.16%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "650.68"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -3.95%  "
$ws.Cells.Item(17, 4).Value = "69.367.67"
$ws.Cells.Item(17, 5).Value = "  -2.13%  "
$ws.Cells.Item(18, 4).Value = "3.472.50"
$ws.Cells.Item(18, 5).Value = "  -4.50%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "12.31"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -5.60%  "
$ws.Cells.Item(20, 5).Value = "  -1.75%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "18.21"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -4.55%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "0.946"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -5.37%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "18.10"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.20%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -2.11%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "98.76"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -6.55%  "
$ws.Cells.Item(26, 5).Value = "  -7.32%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -4.13%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "10.01"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -4.22%  "
$ws.Cells.Item(29, 5).Value = "  -5.04%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "32.53"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -4.30%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "4.25"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -8.73%  "
$ws.Cells.Item(32, 5).Value = "  -6.59%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "11.60"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -5.08%  "
$ws.Cells.Item(34, 5).Value = "  -5.39%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "60.89"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -4.04%  "
$ws.Cells.Item(36, 4).Value = "3.712.35"
$ws.Cells.Item(36, 5).Value = "  -6.58%  "
$ws.Cells.Item(37, 2).Value = "Dai"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.17%  "
$ws.Cells.Item(38, 2).Value = "Bittensor"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "517.85"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +2.40%  "
$ws.Cells.Item(39, 4).Value = "0.0₃0789"
$ws.Cells.Item(39, 5).Value = "  -9.00%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "2.93"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -3.60%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "3.50"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.90%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.374"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -3.91%  "
$ws.Cells.Item(43, 5).Value = "  -2.65%  "
$ws.Cells.Item(44, 2).Value = "CoreDAO"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "3.54"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +72.15%  "
$ws.Cells.Item(45, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "34.28"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -6.88%  "
$ws.Cells.Item(46, 5).Value = "  -3.91%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "3.35"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -4.60%  "
$ws.Cells.Item(48, 5).Value = "  -8.29%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.135"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -4.46%  "
$ws.Cells.Item(50, 5).Value = "  -0.35%  "
$ws.Cells.Item(51, 5).Value = "  -6.08%  "
